$wb = $excel.ActiveWorkbook

# --- "instances" sheet: align gen2 security-group / NIC assignments with gen1 ---
$instances = $wb.Worksheets.Item("instances")
$instances.Range("I3").Value = "sg2"
$instances.Range("I4").Value = "sg1"
$instances.Range("K2").Value = "vsi1nic1:subnet2:sg1:fip2;vsi1nic2:subnet2:sg2"

# --- "menus" sheet: re-sort the "Image Profiles" list (E2:E17) into size order ---
$menus = $wb.Worksheets.Item("menus")
$profiles = @(
    "bx2-2x8",
    "bx2-4x16",
    "bx2-8x32",
    "bx2-16x64",
    "bx2-32x128",
    "bx2-48x192",
    "cx2-2x4",
    "cx2-4x8",
    "cx2-8x16",
    "cx2-16x32",
    "cx2-32x64",
    "mx2-2x16",
    "mx2-4x32",
    "mx2-8x64 " + [char]0x00A0,
    "mx2-16x128",
    "mx2-32x256"
)
for ($i = 0; $i -lt $profiles.Length; $i++) {
    $row = $i + 2
    $menus.Range("E$row").Value = $profiles[$i]
}

$menus.Range("E2").Select()

# Return focus to the originally active "vpc" sheet.
$vpc = $wb.Worksheets.Item("vpc")
$vpc.Activate()
